$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(64,1).Value = "artificial_regional_rare_celltype_diverse"
